$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.272.43'
$ws.Range("E2").Value = '  +0.00%  '

$ws.Range("D3").Value = '3.373.77'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''573.86'
$ws.Range("E5").Value = '  +0.47%  '

$ws.Range("D6").Value = '''136.66'
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").Value = '3.373.84'
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("E9").Value = '  -0.66%  '

$ws.Range("D10").Value = '''7.47'
$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("E12").Value = '  -0.94%  '

$ws.Range("D13").Value = '3.948.56'
$ws.Range("E13").Value = '  +0.19%  '

$ws.Range("D14").Value = '''0.125'
$ws.Range("E14").Value = '  +2.45%  '

$ws.Range("E15").Value = '  +1.59%  '

$ws.Range("D16").Value = '''26.01'
$ws.Range("E16").Value = '  +3.06%  '

$ws.Range("D17").Value = '3.373.10'
$ws.Range("E17").Value = '  +0.22%  '

$ws.Range("D18").Value = '61.382.93'
$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("E19").Value = '  +0.41%  '

$ws.Range("E20").Value = '  +0.88%  '

$ws.Range("E21").Value = '  -1.42%  '

$ws.Range("D22").Value = '''375.57'
$ws.Range("E22").Value = '  -1.68%  '

$ws.Range("E23").Value = '  -3.66%  '

$ws.Range("D24").Value = '3.507.16'
$ws.Range("E24").Value = '  -0.01%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("D26").Value = '''0.0000127'
$ws.Range("E26").Value = '  +7.74%  '

$ws.Range("D27").Value = '''71.49'
$ws.Range("E27").Value = '  +0.91%  '

$ws.Range("E28").Value = '  +2.92%  '

$ws.Range("D29").Value = '''7.49'
$ws.Range("E29").Value = '  -3.84%  '

$ws.Range("E30").Value = '  +0.19%  '

$ws.Range("D31").Value = '''8.25'
$ws.Range("E31").Value = '  +1.38%  '

$ws.Range("E32").Value = '  +3.35%  '

$ws.Range("E33").Value = '  +1.40%  '

$ws.Range("D35").Value = '''23.58'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("D36").Value = '''5.28'
$ws.Range("E36").Value = '  -5.16%  '

$ws.Range("E37").Value = '  -1.90%  '

$ws.Range("D38").Value = '''1.54'
$ws.Range("E38").Value = '  -0.78%  '

$ws.Range("D39").Value = '''165.58'
$ws.Range("E39").Value = '  +1.09%  '

$ws.Range("E40").Value = '  -3.58%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").Value = '''0.774'
$ws.Range("E42").Value = '  +1.72%  '

$ws.Range("E43").Value = '  +5.09%  '

$ws.Range("E44").Value = '  -0.22%  '

$ws.Range("D45").Value = '''41.44'
$ws.Range("E45").Value = '  -0.26%  '

$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").Value = '''24.46'
$ws.Range("E47").Value = '  +5.24%  '

$ws.Range("E48").Value = '  -2.14%  '

$ws.Range("D49").Value = '''22.60'
$ws.Range("E49").Value = '  -2.48%  '

$ws.Range("D50").Value = '2.349.86'
$ws.Range("E50").Value = '  +3.35%  '

$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''2.38'
$ws.Range("E51").Value = '  -1.92%  '
